$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.Value = "'30.022.08"
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -0.90%  '
$cell = $ws.Range('D3')
$cell.Value = "'1.897.01"
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -1.72%  '
$cell = $ws.Range('D4')
$cell.Value = "'0.9987"
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$cell = $ws.Range('D5')
$cell.Value = "'0.7378"
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -2.41%  '
$cell = $ws.Range('D6')
$cell.Value = "'242.85"
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.78%  '
$cell = $ws.Range('D7')
$cell.Value = "'0.9988"
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  -2.54%  '
$cell = $ws.Range('D9')
$cell.Value = "'26.50"
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -3.81%  '
$cell = $ws.Range('D10')
$cell.Value = "'0.06910"
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -1.13%  '
$cell = $ws.Range('D11')
$cell.Value = "'0.7728"
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -1.05%  '
$cell = $ws.Range('D12')
$cell.Value = "'0.07951"
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -0.42%  '
$cell = $ws.Range('D13')
$cell.Value = "'1.902.36"
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -1.44%  '
$cell = $ws.Range('D14')
$cell.Value = "'5.236"
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -2.19%  '
$cell = $ws.Range('D15')
$cell.Value = "'91.64"
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -2.83%  '
$cell = $ws.Range('D16')
$cell.Value = "'30.027.20"
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -0.90%  '
$cell = $ws.Range('D17')
$cell.Value = "'14.17"
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -1.74%  '
$cell = $ws.Range('D18')
$cell.Value = "'5.815"
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range('D19')
$cell.Value = "'0.000007790"
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range('D20')
$cell.Value = "'239.99"
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -4.93%  '
$cell = $ws.Range('D21')
$cell.Value = "'0.9985"
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.05%  '
$cell = $ws.Range('D22')
$cell.Value = "'2.142.77"
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -1.93%  '
$cell = $ws.Range('D23')
$cell.Value = "'0.9984"
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -0.13%  '
$cell = $ws.Range('D24')
$cell.Value = "'6.905"
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +3.45%  '
$cell = $ws.Range('D25')
$cell.Value = "'9.318"
$cell.Style = 'Normal'
$cell = $ws.Range('D26')
$cell.Value = "'167.18"
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +0.82%  '
$cell = $ws.Range('D27')
$cell.Value = "'18.84"
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -0.44%  '
$cell = $ws.Range('D28')
$cell.Value = "'0.1280"
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -3.81%  '
$cell = $ws.Range('D29')
$cell.Value = "'2.037"
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -7.89%  '
$cell = $ws.Range('D30')
$cell.Value = "'1.349"
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -1.07%  '
$cell = $ws.Range('D31')
$cell.Value = "'1.537"
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +1.69%  '
$cell = $ws.Range('D32')
$cell.Value = "'4.303"
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -1.56%  '
$cell = $ws.Range('D33')
$cell.Value = "'4.066"
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('E34').Value = '  -0.85%  '
$cell = $ws.Range('D35')
$cell.Value = "'1.282"
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +0.62%  '
$cell = $ws.Range('D36')
$cell.Value = "'0.7364"
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -1.10%  '
$cell = $ws.Range('D37')
$cell.Value = "'2.712"
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('E38').Value = '  -1.06%  '
$cell = $ws.Range('D39')
$cell.Value = "'2.793"
$cell.Style = 'Normal'
$cell = $ws.Range('D40')
$cell.Value = "'6.317"
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -1.35%  '
$cell = $ws.Range('D41')
$cell.Value = "'74.37"
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -4.45%  '
$cell = $ws.Range('D42')
$cell.Value = "'0.4463"
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -0.06%  '
$cell = $ws.Range('D43')
$cell.Value = "'1.942"
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -1.19%  '
$cell = $ws.Range('D44')
$cell.Value = "'0.9991"
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  +0.54%  '
$cell = $ws.Range('D46')
$cell.Value = "'7.715"
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +3.62%  '
$cell = $ws.Range('D47')
$cell.Value = "'101.26"
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +0.45%  '
$cell = $ws.Range('D48')
$cell.Value = "'9.864"
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +1.37%  '
$cell = $ws.Range('D49')
$cell.Value = "'2.048.34"
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -1.95%  '
$cell = $ws.Range('D50')
$cell.Value = "'36.64"
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -1.41%  '
$cell = $ws.Range('D51')
$cell.Value = "'934.03"
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -5.24%  '
